# Weekly update: insert two new "Pimiento" price records (Zafiro rojo / Zafiro
# verde) right after the existing row 216, pushing the rest of the dataset
# down by two rows (252 data rows -> 254 data rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 217 (row 216 and everything above stays put,
# the previous rows 217.. shift down to 219..).
$ws.Rows.Item(217).Insert()
$ws.Rows.Item(217).Insert()

# ---- New row 217: Zafiro rojo ----
$ws.Cells.Item(217, 1).Value = 11
$ws.Cells.Item(217, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(217, 3).Value = "Bíobío"
$ws.Cells.Item(217, 4).Value = 44617
$ws.Cells.Item(217, 5).Value = 8
$ws.Cells.Item(217, 6).Value = 100112002
$ws.Cells.Item(217, 7).Value = "Pimiento"
$ws.Cells.Item(217, 8).Value = "Zafiro rojo"
$ws.Cells.Item(217, 9).Value = "Primera"
$ws.Cells.Item(217, 10).Value = 100
$ws.Cells.Item(217, 11).Value = 16000
$ws.Cells.Item(217, 12).Value = 17000
$ws.Cells.Item(217, 13).Value = 16500
$ws.Cells.Item(217, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(217, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(217, 16).Value = 1100
$ws.Cells.Item(217, 17).Value = 15
$ws.Cells.Item(217, 18).Value = "Hortaliza"

# ---- New row 218: Zafiro verde ----
$ws.Cells.Item(218, 1).Value = 11
$ws.Cells.Item(218, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(218, 3).Value = "Bíobío"
$ws.Cells.Item(218, 4).Value = 44617
$ws.Cells.Item(218, 5).Value = 8
$ws.Cells.Item(218, 6).Value = 100112002
$ws.Cells.Item(218, 7).Value = "Pimiento"
$ws.Cells.Item(218, 8).Value = "Zafiro verde"
$ws.Cells.Item(218, 9).Value = "Primera"
$ws.Cells.Item(218, 10).Value = 100
$ws.Cells.Item(218, 11).Value = 11000
$ws.Cells.Item(218, 12).Value = 12000
$ws.Cells.Item(218, 13).Value = 11500
$ws.Cells.Item(218, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(218, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(218, 16).Value = 767
$ws.Cells.Item(218, 17).Value = 15
$ws.Cells.Item(218, 18).Value = "Hortaliza"
